# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns with
# the latest scraped values. NumberFormat is forced to Text ("@") before
# assigning D-column values that look numeric so Excel keeps them as the
# original plain-text price strings (e.g. "595.81") instead of silently
# converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.734.18"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.622.59"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.81"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.58"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.69"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.77"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "3.090.77"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "63.563.82"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  +5.28%  "
$ws.Range("D17").Value = "2.614.91"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.35"
$ws.Range("E18").Value = "  +7.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.67"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.51"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.39"
$ws.Range("E25").Value = "  +11.74%  "
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "567.06"
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.24"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.85"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.41"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.12"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.90"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.93"
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("E45").Value = "  +4.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.62"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.631"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0965"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.39"
$ws.Range("E51").Value = "  +3.34%  "
